$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates scraped on 2023-01-03 ("Updated symbol list" GH Action run).
# Columns D (Price) and E (Volume 1h) hold numeric-looking text ("245.25", "-0.54%");
# NumberFormat is forced to "@" (Text) before the assignment so the engine keeps
# them as literal strings instead of auto-coercing to Number/Percentage, matching
# the workbook's original inline-string cells. Columns B/C are plain text and need no such guard.
$updates = @(
    @{ Cell = 'D2'; Value = '245.25'; ForceText = $true }
    @{ Cell = 'E2'; Value = '-0.54%'; ForceText = $true }
    @{ Cell = 'E3'; Value = '-0.06%'; ForceText = $true }
    @{ Cell = 'D4'; Value = '5.271'; ForceText = $true }
    @{ Cell = 'E4'; Value = '1.20%'; ForceText = $true }
    @{ Cell = 'D5'; Value = '0.05722'; ForceText = $true }
    @{ Cell = 'E5'; Value = '0.20%'; ForceText = $true }
    @{ Cell = 'D6'; Value = '6.617'; ForceText = $true }
    @{ Cell = 'E6'; Value = '0.74%'; ForceText = $true }
    @{ Cell = 'D7'; Value = '3.191'; ForceText = $true }
    @{ Cell = 'E7'; Value = '4.32%'; ForceText = $true }
    @{ Cell = 'D8'; Value = '0.8556'; ForceText = $true }
    @{ Cell = 'E8'; Value = '-0.30%'; ForceText = $true }
    @{ Cell = 'D9'; Value = '0.8537'; ForceText = $true }
    @{ Cell = 'E9'; Value = '-2.89%'; ForceText = $true }
    @{ Cell = 'B10'; Value = 'One'; ForceText = $false }
    @{ Cell = 'C10'; Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'; ForceText = $false }
    @{ Cell = 'D10'; Value = '0.01011'; ForceText = $true }
    @{ Cell = 'E10'; Value = '1,582.61%'; ForceText = $true }
    @{ Cell = 'B11'; Value = 'WazirX'; ForceText = $false }
    @{ Cell = 'C11'; Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'; ForceText = $false }
    @{ Cell = 'D11'; Value = '0.1370'; ForceText = $true }
    @{ Cell = 'E11'; Value = '0.24%'; ForceText = $true }
    @{ Cell = 'B12'; Value = 'MandalaExchangeToken'; ForceText = $false }
    @{ Cell = 'C12'; Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'; ForceText = $false }
    @{ Cell = 'D12'; Value = '0.07089'; ForceText = $true }
    @{ Cell = 'E12'; Value = '0.06%'; ForceText = $true }
    @{ Cell = 'B13'; Value = 'BitrueCoin'; ForceText = $false }
    @{ Cell = 'C13'; Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'; ForceText = $false }
    @{ Cell = 'D13'; Value = '0.03157'; ForceText = $true }
    @{ Cell = 'E13'; Value = '9.94%'; ForceText = $true }
    @{ Cell = 'B14'; Value = 'BitMartToken'; ForceText = $false }
    @{ Cell = 'C14'; Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'; ForceText = $false }
    @{ Cell = 'D14'; Value = '0.09323'; ForceText = $true }
    @{ Cell = 'E14'; Value = '-0.66%'; ForceText = $true }
    @{ Cell = 'B15'; Value = 'BitForexToken'; ForceText = $false }
    @{ Cell = 'C15'; Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'; ForceText = $false }
    @{ Cell = 'D15'; Value = '0.001532'; ForceText = $true }
    @{ Cell = 'E15'; Value = '-0.79%'; ForceText = $true }
    @{ Cell = 'D16'; Value = '0.006159'; ForceText = $true }
    @{ Cell = 'E16'; Value = '0.89%'; ForceText = $true }
    @{ Cell = 'D17'; Value = '3.513'; ForceText = $true }
    @{ Cell = 'E17'; Value = '0.80%'; ForceText = $true }
    @{ Cell = 'D18'; Value = '2.182'; ForceText = $true }
    @{ Cell = 'E18'; Value = '-3.83%'; ForceText = $true }
    @{ Cell = 'D19'; Value = '0.3158'; ForceText = $true }
    @{ Cell = 'E19'; Value = '-0.63%'; ForceText = $true }
    @{ Cell = 'D20'; Value = '0.03333'; ForceText = $true }
    @{ Cell = 'E20'; Value = '0.90%'; ForceText = $true }
    @{ Cell = 'D21'; Value = '0.1278'; ForceText = $true }
    @{ Cell = 'E21'; Value = '-1.73%'; ForceText = $true }
    @{ Cell = 'D22'; Value = '3.494'; ForceText = $true }
    @{ Cell = 'E22'; Value = '0.75%'; ForceText = $true }
    @{ Cell = 'D23'; Value = '0.04145'; ForceText = $true }
    @{ Cell = 'E23'; Value = '-0.34%'; ForceText = $true }
    @{ Cell = 'D24'; Value = '0.1363'; ForceText = $true }
    @{ Cell = 'E24'; Value = '-1.21%'; ForceText = $true }
    @{ Cell = 'D25'; Value = '0.001222'; ForceText = $true }
    @{ Cell = 'E25'; Value = '0.01%'; ForceText = $true }
    @{ Cell = 'D26'; Value = '0.004142'; ForceText = $true }
    @{ Cell = 'E26'; Value = '-18.05%'; ForceText = $true }
    @{ Cell = 'E27'; Value = '-0.68%'; ForceText = $true }
    @{ Cell = 'D28'; Value = '0.0001451'; ForceText = $true }
    @{ Cell = 'E28'; Value = '-25.15%'; ForceText = $true }
    @{ Cell = 'D40'; Value = '0.03770'; ForceText = $true }
    @{ Cell = 'E40'; Value = '0.72%'; ForceText = $true }
    @{ Cell = 'D41'; Value = '0.1064'; ForceText = $true }
    @{ Cell = 'E41'; Value = '-0.82%'; ForceText = $true }
    @{ Cell = 'B42'; Value = 'CEJI'; ForceText = $false }
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'; ForceText = $false }
    @{ Cell = 'D42'; Value = '0.002420'; ForceText = $true }
    @{ Cell = 'E42'; Value = '15.25%'; ForceText = $true }
    @{ Cell = 'B43'; Value = 'KickToken'; ForceText = $false }
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'; ForceText = $false }
    @{ Cell = 'D43'; Value = '0.002954'; ForceText = $true }
    @{ Cell = 'E43'; Value = '-47.89%'; ForceText = $true }
    @{ Cell = 'D44'; Value = '0.009431'; ForceText = $true }
    @{ Cell = 'E44'; Value = '0.36%'; ForceText = $true }
    @{ Cell = 'D45'; Value = '0.00005280'; ForceText = $true }
    @{ Cell = 'E45'; Value = '3.54%'; ForceText = $true }
    @{ Cell = 'E46'; Value = '0.17%'; ForceText = $true }
    @{ Cell = 'D47'; Value = '0.09001'; ForceText = $true }
    @{ Cell = 'E47'; Value = '26.82%'; ForceText = $true }
    @{ Cell = 'D48'; Value = '0.002208'; ForceText = $true }
    @{ Cell = 'E48'; Value = '-14.91%'; ForceText = $true }
    @{ Cell = 'D49'; Value = '0.00002103'; ForceText = $true }
    @{ Cell = 'E49'; Value = '0.17%'; ForceText = $true }
    @{ Cell = 'D50'; Value = '0.0002003'; ForceText = $true }
    @{ Cell = 'E50'; Value = '0.17%'; ForceText = $true }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $u.Value
}
